# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list (E16:E22) is re-sorted from descending
# (2403..2309) to ascending (2309..2403) order, and the one-off "Valor
# Mora" amount of 55835 (all the other periods carry 88160) follows the
# 2403 period, which is now the last row instead of the first.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periodos = @("2309", "2310", "2311", "2312", "2401", "2402", "2403")
$valores  = @(88160, 88160, 88160, 88160, 88160, 88160, 55835)

$startRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = $valores[$i]
}

$wb.Save()
